$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.400771666666667
$ws.Range("H2").Value = 19.202315
$ws.Range("I2").Value = 0.4226371084084476
$ws.Range("J2").Value = 0.4226371084084476
$ws.Range("M2").Value = 568.5612486666666
$ws.Range("N2").Value = 1705.683746
$ws.Range("O2").Value = 0.6737621253161296
$ws.Range("P2").Value = 0.6737621253161296
$ws.Range("Q2").Value = 3639.230731230221
$ws.Range("R2").Value = 32753.07658107199
$ws.Range("S2").Value = 0.2847568763987391
$ws.Range("T2").Value = 0.2847568763987391
$ws.Range("G3").Value = 6.400771666666667
$ws.Range("H3").Value = 19.202315
$ws.Range("I3").Value = 0.4226371084084476
$ws.Range("J3").Value = 0.4226371084084476
$ws.Range("M3").Value = 88.00803400000001
$ws.Range("O3").Value = 0.1042921587987053
$ws.Range("P3").Value = 0.1042921587987053
$ws.Range("Q3").Value = 563.3193304662368
$ws.Range("R3").Value = 5069.873974196131
$ws.Range("S3").Value = 0.04407773642435947
$ws.Range("T3").Value = 0.04407773642435947
$ws.Range("G4").Value = 6.400771666666667
$ws.Range("H4").Value = 19.202315
$ws.Range("I4").Value = 0.4226371084084476
$ws.Range("J4").Value = 0.4226371084084476
$ws.Range("M4").Value = 187.2912243333334
$ws.Range("N4").Value = 561.8736730000001
$ws.Range("O4").Value = 0.2219457158851651
$ws.Range("P4").Value = 0.2219457158851651
$ws.Range("Q4").Value = 1198.808362128111
$ws.Range("R4").Value = 10789.275259153
$ws.Range("S4").Value = 0.09380249558534905
$ws.Range("T4").Value = 0.09380249558534903
$ws.Range("I5").Value = 0.3497933355610079
$ws.Range("J5").Value = 0.3497933355610079
$ws.Range("M5").Value = 568.5612486666666
$ws.Range("N5").Value = 1705.683746
$ws.Range("O5").Value = 0.6737621253161296
$ws.Range("P5").Value = 0.6737621253161296
$ws.Range("Q5").Value = 3011.989792251998
$ws.Range("R5").Value = 27107.90813026798
$ws.Range("S5").Value = 0.2356775011890028
$ws.Range("T5").Value = 0.2356775011890028
$ws.Range("I6").Value = 0.3497933355610079
$ws.Range("J6").Value = 0.3497933355610079
$ws.Range("M6").Value = 88.00803400000001
$ws.Range("O6").Value = 0.1042921587987053
$ws.Range("P6").Value = 0.1042921587987053
$ws.Range("Q6").Value = 466.2282219651874
$ws.Range("R6").Value = 4196.053997686687
$ws.Range("S6").Value = 0.03648070209905747
$ws.Range("T6").Value = 0.03648070209905747
$ws.Range("I7").Value = 0.3497933355610079
$ws.Range("J7").Value = 0.3497933355610079
$ws.Range("M7").Value = 187.2912243333334
$ws.Range("N7").Value = 561.8736730000001
$ws.Range("O7").Value = 0.2219457158851651
$ws.Range("P7").Value = 0.2219457158851651
$ws.Range("Q7").Value = 992.187309974599
$ws.Range("R7").Value = 8929.685789771391
$ws.Range("S7").Value = 0.07763513227294769
$ws.Range("T7").Value = 0.07763513227294767
$ws.Range("G8").Value = 3.446504666666667
$ws.Range("H8").Value = 10.339514
$ws.Range("I8").Value = 0.2275695560305444
$ws.Range("J8").Value = 0.2275695560305443
$ws.Range("M8").Value = 568.5612486666666
$ws.Range("N8").Value = 1705.683746
$ws.Range("O8").Value = 0.6737621253161296
$ws.Range("P8").Value = 0.6737621253161296
$ws.Range("Q8").Value = 1959.548996815494
$ws.Range("R8").Value = 17635.94097133945
$ws.Range("S8").Value = 0.1533277477283876
$ws.Range("T8").Value = 0.1533277477283876
$ws.Range("G9").Value = 3.446504666666667
$ws.Range("H9").Value = 10.339514
$ws.Range("I9").Value = 0.2275695560305444
$ws.Range("J9").Value = 0.2275695560305443
$ws.Range("M9").Value = 88.00803400000001
$ws.Range("O9").Value = 0.1042921587987053
$ws.Range("P9").Value = 0.1042921587987053
$ws.Range("Q9").Value = 303.3200998851588
$ws.Range("R9").Value = 2729.880898966429
$ws.Range("S9").Value = 0.02373372027528841
$ws.Range("T9").Value = 0.0237337202752884
$ws.Range("G10").Value = 3.446504666666667
$ws.Range("H10").Value = 10.339514
$ws.Range("I10").Value = 0.2275695560305444
$ws.Range("J10").Value = 0.2275695560305443
$ws.Range("M10").Value = 187.2912243333334
$ws.Range("N10").Value = 561.8736730000001
$ws.Range("O10").Value = 0.2219457158851651
$ws.Range("P10").Value = 0.2219457158851651
$ws.Range("Q10").Value = 645.5000786905471
$ws.Range("R10").Value = 5809.500708214923
$ws.Range("S10").Value = 0.05050808802686835
$ws.Range("T10").Value = 0.05050808802686834
